$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(1,0,0,1,1,2,3,1,2,2,1,0,1,1,1,1,2,0,0,2,1,0,1,2,1,1,1,2,0,0,0,2,0,1,3,1,0,3,1,1,1,1,2,2,2,0,2,1,3,1,1,2,2,5,1,2,3,2,0)

$row = 2
foreach ($v in $values) {
    $ws.Cells.Item($row, 7).Value = $v
    $row++
}
